$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.409.76'
$ws.Range("E2").Value = '  -3.46%  '
$ws.Range("D3").Value = '1.748.95'
$ws.Range("E3").Value = '  -3.72%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4241'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.97%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3598'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07480'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.38%  '
$ws.Range("E10").Value = '  -6.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.098'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.67'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.026'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.222'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.95%  '
$ws.Range("D16").Value = '1.748.81'
$ws.Range("E16").Value = '  -5.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.42'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001064'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06373'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.889'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.88%  '
$ws.Range("D23").Value = '27.463.00'
$ws.Range("E23").Value = '  -3.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.23'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.088'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.73'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.24'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.84%  '
$ws.Range("D28").Value = '1.941.68'
$ws.Range("E28").Value = '  -5.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.120'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.99'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.102'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -9.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.649'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.544'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08869'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.21'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02283'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.20%  '
$ws.Range("E37").Value = '  -4.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05995'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6322'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.945'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.184'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.907'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.387'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5886'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.11%  '
$ws.Range("E47").Value = '  -2.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.89'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.967'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.165'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06817'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.62%  '
